$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1040711.1
$ws.Range("I2").Value = 1818494.8
$ws.Range("K2").Value = 1818494.8
$ws.Range("M2").Value = -1818381.8

$ws.Range("H17").Value = 10795.454
$ws.Range("J17").Value = 11675
$ws.Range("L17").Value = 35025
$ws.Range("N17").Value = -35361

$ws.Range("H33").Value = 7423.125
$ws.Range("J33").Value = 4999.6665
$ws.Range("L33").Value = 4999.6665
$ws.Range("N33").Value = -5457.6665

$ws.Range("H38").Value = 78
$ws.Range("I38").Value = 78
$ws.Range("K38").Value = 234
$ws.Range("M38").Value = 138

$ws.Range("H40").Value = 3499.4614
$ws.Range("I40").Value = 2379.8
$ws.Range("J40").Value = 4199.25
$ws.Range("K40").Value = 2379.8
$ws.Range("L40").Value = 4199.25
$ws.Range("M40").Value = -2204.8
$ws.Range("N40").Value = -4549.25

$ws.Range("H70").Value = 2298.7
$ws.Range("J70").Value = 3099.4
$ws.Range("L70").Value = 9298.200000000001
$ws.Range("N70").Value = -9838.200000000001

$ws.Range("H73").Value = 2298.7
$ws.Range("J73").Value = 3099.4
$ws.Range("L73").Value = 9298.200000000001
$ws.Range("N73").Value = -11170.2

$ws.Range("H76").Value = 4748.5557
$ws.Range("J76").Value = 4949.2
$ws.Range("L76").Value = 4949.2
$ws.Range("N76").Value = -5579.2

$ws.Range("H79").Value = 4748.5557
$ws.Range("J79").Value = 4949.2
$ws.Range("L79").Value = 4949.2
$ws.Range("N79").Value = -7133.2

$ws.Range("H99").Value = 2178.75
$ws.Range("J99").Value = 4935.4
$ws.Range("L99").Value = 14806.2
$ws.Range("N99").Value = -17802.2

$ws.Range("H100").Value = 2928.2083
$ws.Range("I100").Value = 2279.8
$ws.Range("K100").Value = 2279.8
$ws.Range("M100").Value = -1738.8

$ws.Range("H106").Value = 12146.857
$ws.Range("I106").Value = 10439
$ws.Range("J106").Value = 12830
$ws.Range("K106").Value = 10439
$ws.Range("L106").Value = 12830
$ws.Range("M106").Value = -9808
$ws.Range("N106").Value = -14092

$ws.Range("H113").Value = 6630.125
$ws.Range("I113").Value = 6461
$ws.Range("K113").Value = 6461
$ws.Range("M113").Value = -3207

$ws.Range("H116").Value = 5343.8667
$ws.Range("I116").Value = 4867.143
$ws.Range("K116").Value = 4867.143
$ws.Range("M116").Value = -1425.143

$ws.Range("H125").Value = 1390.8572
$ws.Range("I125").Value = 950
$ws.Range("K125").Value = 8550
$ws.Range("M125").Value = -6090

$ws.Range("H135").Value = 2152.7273
$ws.Range("J135").Value = 1779
$ws.Range("L135").Value = 16011
$ws.Range("N135").Value = -21081

$ws.Range("H138").Value = 3236.4614
$ws.Range("I138").Value = 2294.7856
$ws.Range("K138").Value = 6884.3568
$ws.Range("M138").Value = -1744.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15917.546
$ws.Range("I32").Value = 15894
$ws.Range("J32").Value = 16498.334
$ws.Range("K32").Value = 15894
$ws.Range("L32").Value = 16498.334
$ws.Range("M32").Value = -15607
$ws.Range("N32").Value = -17072.334

$ws.Range("H110").Value = 1462.0322
$ws.Range("J110").Value = 423.33334
$ws.Range("L110").Value = 423.33334
$ws.Range("N110").Value = -4513.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 3419.5557
$ws.Range("I94").Value = 3316
$ws.Range("K94").Value = 3316
$ws.Range("M94").Value = -2865

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 594.5
$ws.Range("J23").Value = 594.5
$ws.Range("L23").Value = 1783.5
$ws.Range("N23").Value = -2253.5

$ws.Range("H68").Value = 6002.5
$ws.Range("I68").Value = 5002
$ws.Range("J68").Value = 7003
$ws.Range("K68").Value = 15006
$ws.Range("L68").Value = 21009
$ws.Range("M68").Value = -14195
$ws.Range("N68").Value = -22631

$ws.Range("H71").Value = 6002.5
$ws.Range("I71").Value = 5002
$ws.Range("J71").Value = 7003
$ws.Range("K71").Value = 45018
$ws.Range("L71").Value = 63027
$ws.Range("M71").Value = -40962
$ws.Range("N71").Value = -71139

$ws.Range("H134").Value = 3930.6667
$ws.Range("I134").Value = 3930.6667
$ws.Range("K134").Value = 11792.0001
$ws.Range("M134").Value = -6722.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 287778
$ws.Range("I113").Value = 202263.8
$ws.Range("J113").Value = 501563.5
$ws.Range("K113").Value = 202263.8
$ws.Range("L113").Value = 501563.5
$ws.Range("M113").Value = -200093.8
$ws.Range("N113").Value = -505903.5

$ws.Range("H122").Value = 1856.6
$ws.Range("I122").Value = 1387.7391
$ws.Range("J122").Value = 3397.1428
$ws.Range("K122").Value = 4163.2173
$ws.Range("L122").Value = 10191.4284
$ws.Range("M122").Value = -1713.2173
$ws.Range("N122").Value = -15091.4284

$ws.Range("H132").Value = 31860.059
$ws.Range("I132").Value = 39408.816
$ws.Range("K132").Value = 118226.448
$ws.Range("M132").Value = -115696.448

$ws.Range("H137").Value = 68621.75
$ws.Range("I137").Value = 50709
$ws.Range("J137").Value = 74592.664
$ws.Range("K137").Value = 50709
$ws.Range("L137").Value = 74592.664
$ws.Range("M137").Value = -45609
$ws.Range("N137").Value = -84792.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 82248.42999999999
$ws.Range("I22").Value = 370494.66
$ws.Range("J22").Value = 3635.818
$ws.Range("K22").Value = 370494.66
$ws.Range("L22").Value = 3635.818
$ws.Range("M22").Value = -370199.66
$ws.Range("N22").Value = -4225.818

$ws.Range("H27").Value = 82248.42999999999
$ws.Range("I27").Value = 370494.66
$ws.Range("J27").Value = 3635.818
$ws.Range("K27").Value = 370494.66
$ws.Range("L27").Value = 3635.818
$ws.Range("M27").Value = -370387.66
$ws.Range("N27").Value = -3849.818

$ws.Range("H132").Value = 45921.57
$ws.Range("I132").Value = 54686.957
$ws.Range("K132").Value = 164060.871
$ws.Range("M132").Value = -161530.871

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 853
$ws.Range("I100").Value = 706.25
$ws.Range("K100").Value = 1412.5
$ws.Range("M100").Value = -871.5

$ws.Range("H113").Value = 1328.421
$ws.Range("I113").Value = 1329.375
$ws.Range("K113").Value = 3988.125
$ws.Range("M113").Value = -1818.125

$ws.Range("H126").Value = 127930.305
$ws.Range("I126").Value = 146349.66
$ws.Range("J126").Value = 5134.6665
$ws.Range("K126").Value = 439048.98
$ws.Range("L126").Value = 15403.9995
$ws.Range("M126").Value = -436578.98
$ws.Range("N126").Value = -20343.9995
